$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column (H1), copying the header formatting
# (bold, centered, bordered) from the neighboring header cell G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for row 2.
$ws.Range("H2").Value = 1
